# Update sheet (tab) name and workbook title reference to the new "through" date
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2021-12-05"

# Row 13: only U13/V13 change
$ws.Range("U13").Value = 196
$ws.Range("V13").Value = 0.0249

# Row 14: label + many cells
$ws.Range("A14").Value = "December (through 12-05)"
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 0.3333
$ws.Range("F14").Value = 15
$ws.Range("G14").Value = 0.0625
$ws.Range("H14").Value = 2
$ws.Range("I14").Value = 17
$ws.Range("J14").Value = 0.1053
$ws.Range("O14").Value = 5
$ws.Range("Q14").Value = 2
$ws.Range("R14").Value = 24
$ws.Range("S14").Value = 0.0769
$ws.Range("U14").Value = 37

# Row 15 (Total row): updated sums
$ws.Range("C15").Value = 260
$ws.Range("D15").Value = 0.1156
$ws.Range("F15").Value = 518
$ws.Range("G15").Value = 0.1054
$ws.Range("H15").Value = 65
$ws.Range("I15").Value = 775
$ws.Range("J15").Value = 0.0774
$ws.Range("O15").Value = 485
$ws.Range("P15").Value = 0.1002
$ws.Range("Q15").Value = 66
$ws.Range("R15").Value = 1224
$ws.Range("S15").Value = 0.0512
$ws.Range("U15").Value = 1581
$ws.Range("V15").Value = 0.0589

# New cell S14 needs the percentage number format (same as other arrest_rate cells)
$ws.Range("S14").NumberFormat = $ws.Range("S13").NumberFormat
